# Update countries & provincias Spain
# Applies updated COVID-19 numbers plus the resulting alphabetical
# re-ordering of a couple of rows (Nepal/Portugal and
# Santa Lucia/Timor Oriental), and refreshes the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Septiembre de 2020 a las 14:03"

# 2) Row 5 - India
$ws.Range("B5").Value = 6078200
$ws.Range("C5").Value = 4852
$ws.Range("D5").Value = 5016520
$ws.Range("E5").Value = 966080
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 26
$ws.Range("H5").Value = 95600

# 3) Row 16 - Iran
$ws.Range("B16").Value = 449960
$ws.Range("C16").Value = 3512
$ws.Range("D16").Value = 376531
$ws.Range("E16").Value = 47650
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 190
$ws.Range("H16").Value = 25779

# 4) Row 40 - Kuwait
$ws.Range("B40").Value = 103981
$ws.Range("C40").Value = 437
$ws.Range("D40").Value = 95511
$ws.Range("E40").Value = 7865
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 605

# 5) Rows 51 & 52 - Nepal and Portugal swap alphabetical order, with new data
$ws.Range("A51").Value = "Nepal"
$ws.Range("B51").Value = 74745
$ws.Range("C51").Value = 1351
$ws.Range("D51").Value = 54640
$ws.Range("E51").Value = 19624
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 4
$ws.Range("H51").Value = 481

$ws.Range("A52").Value = "Portugal"
$ws.Range("B52").Value = 73604
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 47647
$ws.Range("E52").Value = 24004
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 1953

# 6) Row 77 - El Salvador
$ws.Range("B77").Value = 28809
$ws.Range("C77").Value = 179
$ws.Range("D77").Value = 23317
$ws.Range("E77").Value = 4661
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 5
$ws.Range("H77").Value = 831

# 7) Row 86 - Republica de Macedonia
$ws.Range("B86").Value = 17674
$ws.Range("C86").Value = 45
$ws.Range("D86").Value = 14642
$ws.Range("E86").Value = 2303
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 4
$ws.Range("H86").Value = 729

# 8) Row 91 - Senegal
$ws.Range("B91").Value = 14919
$ws.Range("C91").Value = 10
$ws.Range("D91").Value = 12231
$ws.Range("E91").Value = 2379
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 309

# 9) Row 113 - Uganda
$ws.Range("B113").Value = 7777
$ws.Range("C113").Value = 247
$ws.Range("D113").Value = 4033
$ws.Range("E113").Value = 3669
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 2
$ws.Range("H113").Value = 75

# 10) Row 145 - Malta
$ws.Range("B145").Value = 3006
$ws.Range("C145").Value = 27
$ws.Range("D145").Value = 2399
$ws.Range("E145").Value = 575
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 32

# 11) Row 168 - Vietnam
$ws.Range("B168").Value = 1077
$ws.Range("C168").Value = 3
$ws.Range("D168").Value = 999
$ws.Range("E168").Value = 43
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 35

# 12) Rows 207 & 208 - Santa Lucia and Timor Oriental swap alphabetical
#     order (values are identical before/after, only the label order
#     changes).
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("B207").Value = 27
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 27
$ws.Range("E207").Value = 0
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0

$ws.Range("A208").Value = "Timor Oriental"
$ws.Range("B208").Value = 27
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 27
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0
